$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price/Volume columns being updated so Excel
# does not coerce numeric-looking strings (e.g. "1.009") into floating point
# numbers or reformat multi-dot values. This preserves the original text-cell
# semantics used throughout the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.003.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.821.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4643'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3634'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07293'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8656'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.882.36'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07569'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.30'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.341'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.464'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008630'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.51%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.258.05'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.147'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.57'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.087.55'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.72'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.857'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.25'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.095'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.31%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.89'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.072'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08904'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.954'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7283'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.137'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.425'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.500'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.074'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05256'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01916'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.930'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.112'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5203'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1633'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.215'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4852'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.61%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.24'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.637'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06219'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.64%  '
